# Adds two new columns, I ("I0") and J ("IF"), to the right of the existing
# H ("IP") column, extending the used range from A1:H57 to A1:J57.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing header cells (B1:H1) use style index 1 (bold font, thin box
# border, centered/top aligned). Copy that formatting from H1 onto the new
# header cells I1:J1 before writing their text, so I1/J1 end up styled the
# same way as the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-57: both I and J get the same value per row.
$values = @(
    6, 8, 9, 7, 9, 9, 7, 8, 8, 8,
    8, 9, 8, 8, 9, 8, 8, 8, 8, 8,
    9, 8, 8, 8, 8, 8, 8, 8, 8, 8,
    9, 7, 8, 10, 7, 7, 8, 8, 8, 8,
    8, 7, 7, 8, 8, 8, 9, 8, 9, 7,
    8, 8, 8, 7, 8, 8
)

for ($idx = 0; $idx -lt $values.Length; $idx++) {
    $row = $idx + 2
    $val = $values[$idx]
    $ws.Cells.Item($row, 9).Value = $val   # column I
    $ws.Cells.Item($row, 10).Value = $val  # column J
}
